$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 25 - this shifts the existing rows 25-33
# down to become rows 26-34 (matching the diff, which is effectively a
# single new weekly record inserted in the middle of the historical series).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new "Región Metropolitana" record.
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 44489
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = 300000000
$ws.Cells.Item(25, 7).Value = "Espárragos"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 350
$ws.Cells.Item(25, 11).Value = 1300
$ws.Cells.Item(25, 12).Value = 1300
$ws.Cells.Item(25, 13).Value = 1300
$ws.Cells.Item(25, 14).Value = "$/kilo"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 1300
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"
